$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Row 2 ("personne_nom"): drop the spellcheck proofErr wrapper, keep the
# existing 3-run split ("p" / "ersonne_" / "nom").
$xml = '<w:p ' + $wns + '><w:r><w:t>p</w:t></w:r><w:r><w:t>ersonne_</w:t></w:r><w:r><w:t>nom</w:t></w:r></w:p>'
$t.Cell(2, 1).Range.Paragraphs.Item(1).Range.InsertXML($xml)

# Row 3 ("personne_prenom"): drop the spellcheck proofErr wrapper, keep the
# existing 2-run split ("personne" / "_prenom").
$xml = '<w:p ' + $wns + '><w:r><w:t>personne</w:t></w:r><w:r><w:t>_prenom</w:t></w:r></w:p>'
$t.Cell(3, 1).Range.Paragraphs.Item(1).Range.InsertXML($xml)

# Row 4 ("personne_naissance"): drop the spellcheck proofErr wrapper, keep
# the existing 2-run split ("personne" / "_naissance").
$xml = '<w:p ' + $wns + '><w:r><w:t>personne</w:t></w:r><w:r><w:t>_naissance</w:t></w:r></w:p>'
$t.Cell(4, 1).Range.Paragraphs.Item(1).Range.InsertXML($xml)

# Row 5: fill in the previously-empty "Rôle_libellé" / "Libellé du rôle" row.
$xml = '<w:p ' + $wns + '><w:r><w:t>Rôle_libe</w:t></w:r><w:r><w:t>llé</w:t></w:r></w:p>'
$t.Cell(5, 1).Range.Paragraphs.Item(1).Range.InsertXML($xml)

$xml = '<w:p ' + $wns + '><w:r><w:t>Libellé du rôle</w:t></w:r></w:p>'
$t.Cell(5, 2).Range.Paragraphs.Item(1).Range.InsertXML($xml)
